# Regenerate the "K" (strikeouts) column (G) for each start row using the
# recomputed values from the box-score reprocessing (K replaces the old
# Strike# derived figure). Only column G, rows 2-68, changes value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2..68 (in row order)
$kValues = @(1,2,1,1,1,0,1,1,0,2,1,0,0,0,2,3,1,1,2,2,1,1,0,1,0,0,1,1,1,3,1,0,2,1,2,1,3,3,1,3,2,1,0,3,1,0,2,1,1,3,1,0,1,0,2,0,0,3,5,2,0,1,3,4,1,2,0)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
